$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation updates to columns E, F, H, I
$ws.Range("F2").Value = 3
$ws.Range("F15").Value = -2
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 8
$ws.Range("F21").Value = 0
$ws.Range("F25").Value = 6
$ws.Range("F27").Value = 2
$ws.Range("F32").Value = 2
$ws.Range("F39").Value = -2
$ws.Range("F40").Value = -2
$ws.Range("F42").Value = 3
$ws.Range("F44").Value = 0
$ws.Range("E50").Value = 1
$ws.Range("H50").Value = -1
$ws.Range("I50").Value = 9
$ws.Range("F53").Value = -1
$ws.Range("F64").Value = -9
$ws.Range("F68").Value = 0
